$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 48
$ws.Cells.Item($row, 1).Value = "III-2021"
$ws.Cells.Item($row, 2).Value = 15952.4
$ws.Cells.Item($row, 3).Value = 9111.1
$ws.Cells.Item($row, 4).Value = 8345.200000000001
$ws.Cells.Item($row, 5).Value = 765.8
$ws.Cells.Item($row, 6).Value = 703.8
$ws.Cells.Item($row, 7).Value = 62.1
$ws.Cells.Item($row, 8).Value = 6841.3
